$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old _GoBack bookmark that sits right after "Week 5".
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. "Continue building core game loop" -> "Expand on Random Map Generation"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Continue building core game loop", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Expand on Random Map Generation", 2) | Out-Null

# ------------------------------------------------------------------
# 3. "Fill in empty rooms with random generation" -> "Introduce chunk map abstraction"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Fill in empty rooms with random generation", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Introduce chunk map abstraction", 2) | Out-Null

# ------------------------------------------------------------------
# 4. "Randomly place key" -> "Separate object generation steps into tiers"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Randomly place key", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Separate object generation steps into tiers", 2) | Out-Null

# ------------------------------------------------------------------
# 5. Delete the "Modify random path generation for multiple paths" bullet entirely.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Modify random path generation for multiple paths*") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 6. Replace the "Extend to 5 level loop" paragraph (which has an embedded
#    proofErr pair splitting it into two runs) with a single clean run
#    reading "Ensure path existence for keys and doors". We do this by
#    inserting a brand-new (proofErr-free) paragraph before it that
#    inherits the same list formatting, filling in the new text, and
#    then deleting the old paragraph outright (which also removes its
#    proofErr markers).
# ------------------------------------------------------------------
$oldPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Extend to 5 level*") {
        $oldPara = $p
        break
    }
}
$origIndex = $oldPara.Index
$insertPoint = $oldPara.Range.Duplicate
$insertPoint.Collapse(1)
$insertPoint.InsertParagraphBefore()

# The freshly-inserted blank paragraph takes the original index; the old
# (proofErr-laden) paragraph is pushed one slot later.
$newParaIndex = $origIndex
$newPara = $d.Paragraphs($newParaIndex)
$newRange = $newPara.Range.Duplicate
$newRange.MoveEnd(1, -1) | Out-Null
$newRange.Text = "Ensure path existence for keys and doors"

$oldParaIndex = $newParaIndex + 1
$oldParaNow = $d.Paragraphs($oldParaIndex)
$oldParaNow.Range.Delete()

# ------------------------------------------------------------------
# 7. "Update Monsters" -> "Finish core game loop"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Update Monsters", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Finish core game loop", 2) | Out-Null

# ------------------------------------------------------------------
# 8. "Add second monster" -> "Set up 5 level stage"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Add second monster", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Set up 5 level stage", 2) | Out-Null

# ------------------------------------------------------------------
# 9. Delete the "Refine monster AI" bullet entirely.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Refine monster AI*") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# Re-add the _GoBack bookmark at its new location: it now starts right
# before "Introduce chunk map abstraction" and ends at the close of the
# "Ensure path existence for keys and doors" bullet.
# ------------------------------------------------------------------
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Introduce chunk map abstraction*") {
        $startPara = $p
    }
    if ($p.Range.Text -like "*Ensure path existence for keys and doors*") {
        $endPara = $p
    }
}
$startRange = $startPara.Range.Duplicate
$startRange.Collapse(1)
$endRange = $endPara.Range.Duplicate
$endRange.MoveEnd(1, -1) | Out-Null
$endRange.Collapse(0)
$bmRange = $d.Range($startRange.Start, $endRange.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 10. Remove the lastRenderedPageBreak marker in front of "Damage application".
# ------------------------------------------------------------------
$d.Content.Find.Execute("Damage application", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Damage application", 2) | Out-Null

Write-Output "edit complete"
